$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1583.8823
$ws.Range("I70").Value = 2800
$ws.Range("J70").Value = 1507.875
$ws.Range("K70").Value = 8400
$ws.Range("L70").Value = 4523.625
$ws.Range("M70").Value = -8130
$ws.Range("N70").Value = -5063.625
$ws.Range("H73").Value = 1583.8823
$ws.Range("I73").Value = 2800
$ws.Range("J73").Value = 1507.875
$ws.Range("K73").Value = 8400
$ws.Range("L73").Value = 4523.625
$ws.Range("M73").Value = -7464
$ws.Range("N73").Value = -6395.625
$ws.Range("H92").Value = 6155.4116
$ws.Range("I92").Value = 6963.467
$ws.Range("J92").Value = 95
$ws.Range("K92").Value = 6963.467
$ws.Range("L92").Value = 95
$ws.Range("M92").Value = -5715.467
$ws.Range("N92").Value = -2591
$ws.Range("H111").Value = 12250
$ws.Range("I111").Value = 14666.667
$ws.Range("J111").Value = 5000
$ws.Range("K111").Value = 44000.001
$ws.Range("L111").Value = 15000
$ws.Range("M111").Value = -40933.001
$ws.Range("N111").Value = -21134
$ws.Range("H112").Value = 1660.6666
$ws.Range("J112").Value = 1863.8695
$ws.Range("L112").Value = 5591.6085
$ws.Range("N112").Value = -7807.6085
$ws.Range("H113").Value = 2877.8125
$ws.Range("J113").Value = 3504.5557
$ws.Range("L113").Value = 3504.5557
$ws.Range("N113").Value = -10012.5557
$ws.Range("H116").Value = 2577.5557
$ws.Range("I116").Value = 2325
$ws.Range("J116").Value = 3082.6667
$ws.Range("K116").Value = 2325
$ws.Range("L116").Value = 3082.6667
$ws.Range("M116").Value = 1117
$ws.Range("N116").Value = -9966.6667
$ws.Range("H132").Value = 1303.8
$ws.Range("I132").Value = 1395.7709
$ws.Range("J132").Value = 935.9167
$ws.Range("K132").Value = 4187.3127
$ws.Range("L132").Value = 2807.7501
$ws.Range("M132").Value = -1657.3127
$ws.Range("N132").Value = -7867.7501
$ws.Range("H137").Value = 2232.238
$ws.Range("I137").Value = 1335.421
$ws.Range("J137").Value = 3595.4
$ws.Range("K137").Value = 4006.263
$ws.Range("L137").Value = 10786.2
$ws.Range("M137").Value = -1456.263
$ws.Range("N137").Value = -15886.2
$ws.Range("H138").Value = 2913.922
$ws.Range("I138").Value = 1377.6875
$ws.Range("J138").Value = 3761.5
$ws.Range("K138").Value = 4133.0625
$ws.Range("L138").Value = 11284.5
$ws.Range("M138").Value = 1006.9375
$ws.Range("N138").Value = -21564.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 455.33334
$ws.Range("I5").Value = 443
$ws.Range("J5").Value = 498.5
$ws.Range("K5").Value = 443
$ws.Range("L5").Value = 498.5
$ws.Range("M5").Value = -331
$ws.Range("N5").Value = -722.5
$ws.Range("H32").Value = 24238.91
$ws.Range("I32").Value = 27158.543
$ws.Range("J32").Value = 10808.6
$ws.Range("K32").Value = 27158.543
$ws.Range("L32").Value = 10808.6
$ws.Range("M32").Value = -26871.543
$ws.Range("N32").Value = -11382.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 455.33334
$ws.Range("I4").Value = 443
$ws.Range("J4").Value = 498.5
$ws.Range("K4").Value = 443
$ws.Range("L4").Value = 498.5
$ws.Range("M4").Value = -328
$ws.Range("N4").Value = -728.5
$ws.Range("H134").Value = 24463.533
$ws.Range("I134").Value = 2082.8064
$ws.Range("J134").Value = 74020.86
$ws.Range("K134").Value = 6248.4192
$ws.Range("L134").Value = 222062.58
$ws.Range("M134").Value = -3713.4192
$ws.Range("N134").Value = -227132.58

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 131.83333
$ws.Range("I7").Value = 128.4
$ws.Range("J7").Value = 134.28572
$ws.Range("K7").Value = 128.4
$ws.Range("L7").Value = 134.28572
$ws.Range("M7").Value = -15.40000000000001
$ws.Range("N7").Value = -360.28572
$ws.Range("H22").Value = 235.22223
$ws.Range("I22").Value = 183.8
$ws.Range("J22").Value = 299.5
$ws.Range("K22").Value = 183.8
$ws.Range("L22").Value = 299.5
$ws.Range("M22").Value = 166.2
$ws.Range("N22").Value = -999.5
$ws.Range("H31").Value = 9305.777
$ws.Range("I31").Value = 11462.818
$ws.Range("J31").Value = 5916.143
$ws.Range("K31").Value = 11462.818
$ws.Range("L31").Value = 5916.143
$ws.Range("M31").Value = -11167.818
$ws.Range("N31").Value = -6506.143
$ws.Range("H34").Value = 9305.777
$ws.Range("I34").Value = 11462.818
$ws.Range("J34").Value = 5916.143
$ws.Range("K34").Value = 11462.818
$ws.Range("L34").Value = 5916.143
$ws.Range("M34").Value = -11260.818
$ws.Range("N34").Value = -6320.143
$ws.Range("H58").Value = 1784861.4
$ws.Range("I58").Value = 2842546
$ws.Range("J58").Value = 3497.8948
$ws.Range("K58").Value = 2842546
$ws.Range("L58").Value = 3497.8948
$ws.Range("M58").Value = -2842343
$ws.Range("N58").Value = -3903.8948
$ws.Range("H62").Value = 2945.4546
$ws.Range("I62").Value = 2988.889
$ws.Range("K62").Value = 2988.889
$ws.Range("M62").Value = -2364.889
$ws.Range("H65").Value = 2945.4546
$ws.Range("I65").Value = 2988.889
$ws.Range("K65").Value = 14944.445
$ws.Range("M65").Value = -11824.445
$ws.Range("H86").Value = 3127.9312
$ws.Range("I86").Value = 3117.8333
$ws.Range("J86").Value = 3144.4546
$ws.Range("K86").Value = 3117.8333
$ws.Range("L86").Value = 3144.4546
$ws.Range("M86").Value = -1994.8333
$ws.Range("N86").Value = -5390.4546
$ws.Range("H89").Value = 3127.9312
$ws.Range("I89").Value = 3117.8333
$ws.Range("J89").Value = 3144.4546
$ws.Range("K89").Value = 15589.1665
$ws.Range("L89").Value = 15722.273
$ws.Range("M89").Value = -9973.166499999999
$ws.Range("N89").Value = -26954.273
$ws.Range("H134").Value = 2534.5645
$ws.Range("I134").Value = 1655.3611
$ws.Range("K134").Value = 4966.0833
$ws.Range("M134").Value = -2431.0833
$ws.Range("H136").Value = 1784861.4
$ws.Range("I136").Value = 2842546
$ws.Range("J136").Value = 3497.8948
$ws.Range("K136").Value = 8527638
$ws.Range("L136").Value = 10493.6844
$ws.Range("M136").Value = -8525088
$ws.Range("N136").Value = -15593.6844
$ws.Range("H140").Value = 36599
$ws.Range("J140").Value = 36599
$ws.Range("L140").Value = 36599
$ws.Range("N140").Value = -46959
$ws.Range("H141").Value = 29664.705
$ws.Range("I141").Value = 20296
$ws.Range("J141").Value = 30250.25
$ws.Range("K141").Value = 20296
$ws.Range("L141").Value = 30250.25
$ws.Range("M141").Value = -15116
$ws.Range("N141").Value = -40610.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3403957.8
$ws.Range("I5").Value = 353.87097
$ws.Range("J5").Value = 9265720
$ws.Range("K5").Value = 1061.61291
$ws.Range("L5").Value = 27797160
$ws.Range("M5").Value = -949.6129100000001
$ws.Range("N5").Value = -27797384
$ws.Range("H131").Value = 44583.91
$ws.Range("I131").Value = 1479.3334
$ws.Range("J131").Value = 96309.39999999999
$ws.Range("K131").Value = 4438.0002
$ws.Range("L131").Value = 288928.2
$ws.Range("M131").Value = 601.9997999999996
$ws.Range("N131").Value = -299008.2
$ws.Range("H135").Value = 3403957.8
$ws.Range("I135").Value = 353.87097
$ws.Range("J135").Value = 9265720
$ws.Range("K135").Value = 3184.83873
$ws.Range("L135").Value = 83391480
$ws.Range("M135").Value = -649.8387299999999
$ws.Range("N135").Value = -83396550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1641.7858
$ws.Range("I97").Value = 1150.9524
$ws.Range("J97").Value = 3114.2856
$ws.Range("K97").Value = 1150.9524
$ws.Range("L97").Value = 3114.2856
$ws.Range("M97").Value = -654.9523999999999
$ws.Range("N97").Value = -4106.2856
$ws.Range("H107").Value = 375.42856
$ws.Range("I107").Value = 174.33333
$ws.Range("J107").Value = 878.1667
$ws.Range("K107").Value = 174.33333
$ws.Range("L107").Value = 878.1667
$ws.Range("M107").Value = 1745.66667
$ws.Range("N107").Value = -4718.1667
$ws.Range("H113").Value = 3411.3
$ws.Range("I113").Value = 10000
$ws.Range("J113").Value = 2679.2222
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 2679.2222
$ws.Range("M113").Value = -7830
$ws.Range("N113").Value = -7019.2222
$ws.Range("H118").Value = 20206.666
$ws.Range("J118").Value = 20206.666
$ws.Range("L118").Value = 20206.666
$ws.Range("N118").Value = -23520.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1266.7858
$ws.Range("I68").Value = 861.1111
$ws.Range("J68").Value = 1997
$ws.Range("K68").Value = 861.1111
$ws.Range("L68").Value = 1997
$ws.Range("M68").Value = -112.1111
$ws.Range("N68").Value = -3495
$ws.Range("H71").Value = 1266.7858
$ws.Range("I71").Value = 861.1111
$ws.Range("J71").Value = 1997
$ws.Range("K71").Value = 4305.555499999999
$ws.Range("L71").Value = 9985
$ws.Range("M71").Value = -561.5554999999995
$ws.Range("N71").Value = -17473

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3473.762
$ws.Range("I62").Value = 3337.375
$ws.Range("J62").Value = 3557.6924
$ws.Range("K62").Value = 3337.375
$ws.Range("L62").Value = 3557.6924
$ws.Range("M62").Value = -2713.375
$ws.Range("N62").Value = -4805.6924
$ws.Range("H65").Value = 3473.762
$ws.Range("I65").Value = 3337.375
$ws.Range("J65").Value = 3557.6924
$ws.Range("K65").Value = 16686.875
$ws.Range("L65").Value = 17788.462
$ws.Range("M65").Value = -13566.875
$ws.Range("N65").Value = -24028.462
$ws.Range("H107").Value = 1463.5
$ws.Range("I107").Value = 605.2308
$ws.Range("J107").Value = 3057.4285
$ws.Range("K107").Value = 1815.6924
$ws.Range("L107").Value = 9172.2855
$ws.Range("M107").Value = 104.3075999999999
$ws.Range("N107").Value = -13012.2855
$ws.Range("H113").Value = 4548.222
$ws.Range("I113").Value = 6017
$ws.Range("J113").Value = 1610.6666
$ws.Range("K113").Value = 18051
$ws.Range("L113").Value = 4831.9998
$ws.Range("M113").Value = -15881
$ws.Range("N113").Value = -9171.9998
$ws.Range("H126").Value = 1653.48
$ws.Range("I126").Value = 1562.625
$ws.Range("J126").Value = 1815
$ws.Range("K126").Value = 4687.875
$ws.Range("L126").Value = 5445
$ws.Range("M126").Value = -2217.875
$ws.Range("N126").Value = -10385
$ws.Range("H132").Value = 1355.2549
$ws.Range("I132").Value = 605.0625
$ws.Range("J132").Value = 2618.7368
$ws.Range("K132").Value = 1815.1875
$ws.Range("L132").Value = 7856.2104
$ws.Range("M132").Value = 714.8125
$ws.Range("N132").Value = -12916.2104
$ws.Range("H136").Value = 3358.805
$ws.Range("I136").Value = 1615.7142
$ws.Range("J136").Value = 7113.154
$ws.Range("K136").Value = 4847.142599999999
$ws.Range("L136").Value = 21339.462
$ws.Range("M136").Value = -2297.142599999999
$ws.Range("N136").Value = -26439.462
